$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.965.02'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '1.628.12'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.24'
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0619'
$ws.Range("E9").Value = '  -3.23%  '
$ws.Range("E10").Value = '  -5.62%  '
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("D12").Value = '1.853.72'
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("D13").Value = '1.627.02'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("E14").Value = '  -2.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.528'
$ws.Range("E15").Value = '  -3.07%  '
$ws.Range("D16").Value = '25.965.49'
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").Value = '0.0₃0739'
$ws.Range("E17").Value = '  -3.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.40'
$ws.Range("E18").Value = '  -3.20%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.00'
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("E21").Value = '  -2.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.57'
$ws.Range("E22").Value = '  -3.60%  '
$ws.Range("E23").Value = '  -2.14%  '
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.21'
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("E27").Value = '  -3.56%  '
$ws.Range("E28").Value = '  -2.03%  '
$ws.Range("E29").Value = '  -1.59%  '
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("E31").Value = '  -2.13%  '
$ws.Range("E32").Value = '  -4.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.11'
$ws.Range("E33").Value = '  -5.42%  '
$ws.Range("E34").Value = '  -2.58%  '
$ws.Range("E35").Value = '  -2.67%  '
$ws.Range("D36").Value = '1.124.63'
$ws.Range("E36").Value = '  -0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.851'
$ws.Range("E37").Value = '  -6.06%  '
$ws.Range("E38").Value = '  -1.59%  '
$ws.Range("E39").Value = '  -3.25%  '
$ws.Range("E40").Value = '  -2.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.10'
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.769'
$ws.Range("E42").Value = '  -3.54%  '
$ws.Range("B43").Value = 'RocketPoolETH'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D43").Value = '1.764.35'
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.16'
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₆0114'
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0533'
$ws.Range("E46").Value = '  +1.96%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.48'
$ws.Range("E47").Value = '  -3.52%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.48'
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.413'
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.52'
$ws.Range("E50").Value = '  -2.93%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.01'
$ws.Range("E51").Value = '  +0.18%  '
